# Generate Report for Handoff
# Inserts a new "bc2f3d67-d464-4b0f-b2a7-526819421d68" entry as the new
# first data row on every worksheet, pushing the existing
# "db9f6251-24b7-4c0c-b15f-c9c0a6d91a79" entry down to row 3.

$wb = $excel.ActiveWorkbook

# Colour used by the workbook's "HyperLink" cell style (cornflower blue,
# 0x6495ED) encoded as a BGR integer the way the Excel object model wants it.
$hyperlinkColor = 15570276

function Style-AsHyperlink($rng) {
    $rng.Font.Color = $hyperlinkColor
    $rng.Font.Underline = 2
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Push the existing row 2 data down to row 3.
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-31-12 20:31:35"

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/fef31b5d44040e396d7f1d6d43cd37b8d06d90b4/e2e/db9f6251-24b7-4c0c-b15f-c9c0a6d91a79.md", "", "", "db9f6251-24b7-4c0c-b15f-c9c0a6d91a79.md")
Style-AsHyperlink $ws.Range("A3")

# Remove the old hyperlink on row 2 so it can be replaced below.
$ws.Range("A2").Hyperlinks.Delete()

# New row 2 data.
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = "2016-31-12 20:31:49"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/fef31b5d44040e396d7f1d6d43cd37b8d06d90b4/e2e/bc2f3d67-d464-4b0f-b2a7-526819421d68.md", "", "", "bc2f3d67-d464-4b0f-b2a7-526819421d68.md")
Style-AsHyperlink $ws.Range("A2")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

# Push the existing row 2 data down to row 3.
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "2016-03-12 20:31:30"
$ws.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("I3").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/fef31b5d44040e396d7f1d6d43cd37b8d06d90b4/e2e/db9f6251-24b7-4c0c-b15f-c9c0a6d91a79.md", "", "", "db9f6251-24b7-4c0c-b15f-c9c0a6d91a79.md")
Style-AsHyperlink $ws.Range("A3")

$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/fef31b5d44040e396d7f1d6d43cd37b8d06d90b4/e2e/db9f6251-24b7-4c0c-b15f-c9c0a6d91a79.md", "", "", ".md")
Style-AsHyperlink $ws.Range("B3")

$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fee19c433ffe34b7bfd79ce13146b46788be0f6a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/db9f6251-24b7-4c0c-b15f-c9c0a6d91a79.ca4012b9ad1747aaae49af8643dbb9373204d0b4.zh-cn.xlf", "", "", "db9f6251-24b7-4c0c-b15f-c9c0a6d91a79.ca4012b9ad1747aaae49af8643dbb9373204d0b4.zh-cn.xlf")
Style-AsHyperlink $ws.Range("D3")

# Remove the old hyperlinks on row 2 so they can be replaced below.
$ws.Range("A2").Hyperlinks.Delete()
$ws.Range("B2").Hyperlinks.Delete()
$ws.Range("D2").Hyperlinks.Delete()

# New row 2 data.
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("E2").Value = "2016-03-12 20:31:46"
$ws.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("H2").Value = "0001-01-01 00:00:00"
$ws.Range("I2").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/fef31b5d44040e396d7f1d6d43cd37b8d06d90b4/e2e/bc2f3d67-d464-4b0f-b2a7-526819421d68.md", "", "", "bc2f3d67-d464-4b0f-b2a7-526819421d68.md")
Style-AsHyperlink $ws.Range("A2")

$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/fef31b5d44040e396d7f1d6d43cd37b8d06d90b4/e2e/bc2f3d67-d464-4b0f-b2a7-526819421d68.md", "", "", ".md")
Style-AsHyperlink $ws.Range("B2")

$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fee19c433ffe34b7bfd79ce13146b46788be0f6a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/bc2f3d67-d464-4b0f-b2a7-526819421d68.85ca65e91160a90df44eac3a42a67869e984e0fb.zh-cn.xlf", "", "", "bc2f3d67-d464-4b0f-b2a7-526819421d68.85ca65e91160a90df44eac3a42a67869e984e0fb.zh-cn.xlf")
Style-AsHyperlink $ws.Range("D2")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

# Push the existing row 2 data down to row 3.
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "2016-03-12 20:31:35"
$ws.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("I3").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/fef31b5d44040e396d7f1d6d43cd37b8d06d90b4/e2e/db9f6251-24b7-4c0c-b15f-c9c0a6d91a79.md", "", "", "db9f6251-24b7-4c0c-b15f-c9c0a6d91a79.md")
Style-AsHyperlink $ws.Range("A3")

$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/fef31b5d44040e396d7f1d6d43cd37b8d06d90b4/e2e/db9f6251-24b7-4c0c-b15f-c9c0a6d91a79.md", "", "", ".md")
Style-AsHyperlink $ws.Range("B3")

$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b07029b0812f569ddfd443b9664291286081f5b6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/db9f6251-24b7-4c0c-b15f-c9c0a6d91a79.ca4012b9ad1747aaae49af8643dbb9373204d0b4.de-de.xlf", "", "", "db9f6251-24b7-4c0c-b15f-c9c0a6d91a79.ca4012b9ad1747aaae49af8643dbb9373204d0b4.de-de.xlf")
Style-AsHyperlink $ws.Range("D3")

# Remove the old hyperlinks on row 2 so they can be replaced below.
$ws.Range("A2").Hyperlinks.Delete()
$ws.Range("B2").Hyperlinks.Delete()
$ws.Range("D2").Hyperlinks.Delete()

# New row 2 data.
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("E2").Value = "2016-03-12 20:31:49"
$ws.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("H2").Value = "0001-01-01 00:00:00"
$ws.Range("I2").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/fef31b5d44040e396d7f1d6d43cd37b8d06d90b4/e2e/bc2f3d67-d464-4b0f-b2a7-526819421d68.md", "", "", "bc2f3d67-d464-4b0f-b2a7-526819421d68.md")
Style-AsHyperlink $ws.Range("A2")

$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/fef31b5d44040e396d7f1d6d43cd37b8d06d90b4/e2e/bc2f3d67-d464-4b0f-b2a7-526819421d68.md", "", "", ".md")
Style-AsHyperlink $ws.Range("B2")

$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b07029b0812f569ddfd443b9664291286081f5b6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/bc2f3d67-d464-4b0f-b2a7-526819421d68.85ca65e91160a90df44eac3a42a67869e984e0fb.de-de.xlf", "", "", "bc2f3d67-d464-4b0f-b2a7-526819421d68.85ca65e91160a90df44eac3a42a67869e984e0fb.de-de.xlf")
Style-AsHyperlink $ws.Range("D2")
